# Update column F (dSF) values on Sheet1 to reflect re-pulled data / new mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of Excel row number -> new value for column F
$updates = @{
    2  = -4
    13 = -7
    19 = -4
    21 = -4
    23 = -3
    25 = -1
    28 = -1
    32 = -4
    34 = -3
    36 = -6
    37 = -2
    42 = -6
    43 = -9
    44 = -2
    46 = -2
    47 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
